# [Typed SDK] Fix problems found during testing Explorer example
#
# 1. Rename the "Sheet1" worksheet to "List of classes".
# 2. Fix a typo in the Events column of the "List of classes" sheet:
#    "Many manye events" -> "Many many events" (cell G12).
# 3. Make "List of classes" the active/selected tab (was "Intro"),
#    with its own selection set to G13.
# 4. Leave "Intro" selection at A63 (as it was), but no longer the
#    active tab.

$wb = $excel.ActiveWorkbook

$wsIntro = $wb.Worksheets.Item("Intro")
$wsList  = $wb.Worksheets.Item("Sheet1")

# Rename "Sheet1" -> "List of classes"
$wsList.Name = "List of classes"

# Fix the typo "Many manye events" -> "Many many events"
$wsList.Range("G12").Value = "Many many events"

# Restore/keep the Intro sheet's own selection (A63) while it is still
# active, then hand the active tab over to "List of classes".
$wsIntro.Activate()
$wsIntro.Range("A63").Select()

$wsList.Activate()
$wsList.Range("G13").Select()
